$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header relabel: insert Harvard case classification column, shifting old "average_doctor" to "_old" ---
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# --- Updated per-row statistics (recomputed after adding Harvard case classification) ---
$ws.Range("E4").Value = 0.405
$ws.Range("F4").Value = 0.076
$ws.Range("G4").Value = 0.276
$ws.Range("N4").Value = 0.405
$ws.Range("O4").Value = 0.062
$ws.Range("P4").Value = 0.25
$ws.Range("W4").Value = 0.242
$ws.Range("X4").Value = 0.107
$ws.Range("Y4").Value = 0.328
$ws.Range("AI4").Value = 0.235
$ws.Range("AK4").Value = 0.26
$ws.Range("AU4").Value = 0.151
$ws.Range("AV4").Value = 0.029
$ws.Range("AW4").Value = 0.17
$ws.Range("BA4").Value = 1.922
$ws.Range("BB4").Value = 0.155
$ws.Range("BC4").Value = 0.393
$ws.Range("BG4").Value = 0.714
$ws.Range("BH4").Value = 0.145
$ws.Range("BI4").Value = 0.38
$ws.Range("BM4").Value = 0.6879999999999999
$ws.Range("BN4").Value = 0.081
$ws.Range("BO4").Value = 0.285
$ws.Range("BP4").Value = 0.641
$ws.Range("BQ4").Value = 0.643
$ws.Range("E5").Value = 0.524
$ws.Range("F5").Value = 0.093
$ws.Range("G5").Value = 0.306
$ws.Range("N5").Value = 0.76
$ws.Range("O5").Value = 0.08
$ws.Range("P5").Value = 0.282
$ws.Range("W5").Value = 0.244
$ws.Range("X5").Value = 0.112
$ws.Range("Y5").Value = 0.335
$ws.Range("AI5").Value = 0.28
$ws.Range("AJ5").Value = 0.098
$ws.Range("AK5").Value = 0.314
$ws.Range("AU5").Value = 0.294
$ws.Range("AV5").Value = 0.095
$ws.Range("AW5").Value = 0.308
$ws.Range("BA5").Value = 1.364
$ws.Range("BB5").Value = 0.09
$ws.Range("BC5").Value = 0.299
$ws.Range("BG5").Value = 0.392
$ws.Range("BH5").Value = 0.048
$ws.Range("BI5").Value = 0.22
$ws.Range("BM5").Value = 0.576
$ws.Range("BN5").Value = 0.076
$ws.Range("BO5").Value = 0.275
$ws.Range("BP5").Value = 0.455
$ws.Range("BQ5").Value = 0.456
$ws.Range("E6").Value = 0.457
$ws.Range("N6").Value = 0.528
$ws.Range("W6").Value = 0.243
$ws.Range("AI6").Value = 0.256
$ws.Range("AU6").Value = 0.2
$ws.Range("BA6").Value = 1.583
$ws.Range("BG6").Value = 0.506
$ws.Range("BM6").Value = 0.627
$ws.Range("BP6").Value = 0.528
$ws.Range("BQ6").Value = 0.53
$ws.Range("E7").Value = 0.495
$ws.Range("N7").Value = 0.647
$ws.Range("W7").Value = 0.244
$ws.Range("AI7").Value = 0.27
$ws.Range("AU7").Value = 0.247
$ws.Range("BA7").Value = 1.442
$ws.Range("BG7").Value = 0.431
$ws.Range("BM7").Value = 0.595
$ws.Range("BP7").Value = 0.481
$ws.Range("BQ7").Value = 0.483
$ws.Range("E8").Value = 0.55
$ws.Range("F8").Value = 0.12
$ws.Range("G8").Value = 0.346
$ws.Range("N8").Value = 0.762
$ws.Range("O8").Value = 0.07000000000000001
$ws.Range("P8").Value = 0.264
$ws.Range("W8").Value = 0.245
$ws.Range("X8").Value = 0.114
$ws.Range("Y8").Value = 0.337
$ws.Range("AI8").Value = 0.256
$ws.Range("AJ8").Value = 0.099
$ws.Range("AK8").Value = 0.315
$ws.Range("AU8").Value = 0.229
$ws.Range("AV8").Value = 0.067
$ws.Range("AW8").Value = 0.259
$ws.Range("BA8").Value = 1.688
$ws.Range("BB8").Value = 0.132
$ws.Range("BC8").Value = 0.363
$ws.Range("BG8").Value = 0.535
$ws.Range("BH8").Value = 0.105
$ws.Range("BI8").Value = 0.324
$ws.Range("BM8").Value = 0.704
$ws.Range("BN8").Value = 0.068
$ws.Range("BO8").Value = 0.261
$ws.Range("BP8").Value = 0.5629999999999999
$ws.Range("BQ8").Value = 0.57
$ws.Range("E9").Value = 0.452
$ws.Range("F9").Value = 0.248
$ws.Range("G9").Value = 0.498
$ws.Range("N9").Value = 0.643
$ws.Range("O9").Value = 0.23
$ws.Range("P9").Value = 0.479
$ws.Range("W9").Value = 0.143
$ws.Range("X9").Value = 0.122
$ws.Range("Y9").Value = 0.35
$ws.Range("AI9").Value = 0.143
$ws.Range("AJ9").Value = 0.122
$ws.Range("AK9").Value = 0.35
$ws.Range("BA9").Value = 1.572
$ws.Range("BB9").Value = 0.241
$ws.Range("BC9").Value = 0.491
$ws.Range("BG9").Value = 0.548
$ws.Range("BM9").Value = 0.619
$ws.Range("BN9").Value = 0.236
$ws.Range("BO9").Value = 0.486
$ws.Range("BP9").Value = 0.524
$ws.Range("BQ9").Value = 0.519
$ws.Range("E10").Value = 0.595
$ws.Range("F10").Value = 0.241
$ws.Range("G10").Value = 0.491
$ws.Range("N10").Value = 0.857
$ws.Range("O10").Value = 0.122
$ws.Range("P10").Value = 0.35
$ws.Range("W10").Value = 0.286
$ws.Range("X10").Value = 0.204
$ws.Range("Y10").Value = 0.452
$ws.Range("AI10").Value = 0.286
$ws.Range("AJ10").Value = 0.204
$ws.Range("AK10").Value = 0.452
$ws.Range("AU10").Value = 0.214
$ws.Range("AV10").Value = 0.168
$ws.Range("AW10").Value = 0.41
$ws.Range("BA10").Value = 1.928
$ws.Range("BB10").Value = 0.249
$ws.Range("BC10").Value = 0.499
$ws.Range("BG10").Value = 0.595
$ws.Range("BH10").Value = 0.241
$ws.Range("BI10").Value = 0.491
$ws.Range("BM10").Value = 0.857
$ws.Range("BN10").Value = 0.122
$ws.Range("BO10").Value = 0.35
$ws.Range("BP10").Value = 0.643
$ws.Range("BQ10").Value = 0.667
$ws.Range("E11").Value = 0.619
$ws.Range("F11").Value = 0.236
$ws.Range("G11").Value = 0.486
$ws.Range("N11").Value = 0.881
$ws.Range("O11").Value = 0.105
$ws.Range("P11").Value = 0.324
$ws.Range("W11").Value = 0.286
$ws.Range("X11").Value = 0.204
$ws.Range("Y11").Value = 0.452
$ws.Range("AI11").Value = 0.286
$ws.Range("AJ11").Value = 0.204
$ws.Range("AK11").Value = 0.452
$ws.Range("AU11").Value = 0.31
$ws.Range("AV11").Value = 0.214
$ws.Range("AW11").Value = 0.462
$ws.Range("BA11").Value = 1.928
$ws.Range("BB11").Value = 0.249
$ws.Range("BC11").Value = 0.499
$ws.Range("BG11").Value = 0.595
$ws.Range("BH11").Value = 0.241
$ws.Range("BI11").Value = 0.491
$ws.Range("BM11").Value = 0.857
$ws.Range("BN11").Value = 0.122
$ws.Range("BO11").Value = 0.35
$ws.Range("BP11").Value = 0.643
$ws.Range("BQ11").Value = 0.667
$ws.Range("E12").Value = 1.462
$ws.Range("F12").Value = 0.71
$ws.Range("G12").Value = 0.843
$ws.Range("N12").Value = 1.526
$ws.Range("O12").Value = 1.091
$ws.Range("P12").Value = 1.045
$ws.Range("W12").Value = 1.75
$ws.Range("X12").Value = 0.6879999999999999
$ws.Range("Y12").Value = 0.829
$ws.Range("AI12").Value = 1.917
$ws.Range("AJ12").Value = 0.91
$ws.Range("AK12").Value = 0.954
$ws.Range("AU12").Value = 2.692
$ws.Range("AV12").Value = 1.598
$ws.Range("AW12").Value = 1.264
$ws.Range("BA12").Value = 3.691
$ws.Range("BB12").Value = 0.388
$ws.Range("BC12").Value = 0.622
$ws.Range("BG12").Value = 1.08
$ws.Range("BH12").Value = 0.074
$ws.Range("BI12").Value = 0.271
$ws.Range("BM12").Value = 1.361
$ws.Range("BN12").Value = 0.397
$ws.Range("BO12").Value = 0.63
$ws.Range("BP12").Value = 1.23
$ws.Range("BQ12").Value = 1.29
$ws.Range("E13").Value = 1.637
$ws.Range("F13").Value = 0.711
$ws.Range("G13").Value = 0.843
$ws.Range("N13").Value = 2.319
$ws.Range("O13").Value = 1.25
$ws.Range("P13").Value = 1.118
$ws.Range("W13").Value = 1.063
$ws.Range("X13").Value = 0.161
$ws.Range("Y13").Value = 0.401
$ws.Range("AI13").Value = 1.365
$ws.Range("AJ13").Value = 0.41
$ws.Range("AK13").Value = 0.64
$ws.Range("AU13").Value = 2.45
$ws.Range("AV13").Value = 1.478
$ws.Range("AW13").Value = 1.216
$ws.Range("BA13").Value = 2.54
$ws.Range("BB13").Value = 0.346
$ws.Range("BC13").Value = 0.589
$ws.Range("BG13").Value = 0.613
$ws.Range("BH13").Value = 0.054
$ws.Range("BI13").Value = 0.232
$ws.Range("BM13").Value = 0.975
$ws.Range("BN13").Value = 0.313
$ws.Range("BO13").Value = 0.5590000000000001
$ws.Range("BP13").Value = 0.847
$ws.Range("BQ13").Value = 0.806
